$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sheet 1")

# Header row
$ws.Range("A1").Value = "UserName"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "UserRole"
$ws.Range("D1").Value = "Email"
$ws.Range("E1").Value = "MobileNo"

# Row 2
$ws.Range("A2").Value = "shubham"
$ws.Range("B2").Value = "dev_123"
$ws.Range("C2").Value = "developers"
$ws.Range("D2").Value = "shubham@blog.com"
$ws.Range("E2").Value = 1122334455

# Row 3
$ws.Range("A3").Value = "gaurav"
$ws.Range("B3").Value = "check_1233"
$ws.Range("C3").Value = "tester"
$ws.Range("D3").Value = "test@test.com"
$ws.Range("E3").Value = 9988273

# Row 4
$ws.Range("A4").Value = "gaurav"
$ws.Range("B4").Value = "check_1233"
$ws.Range("C4").Value = "tester"
$ws.Range("D4").Value = "test@test.com"
$ws.Range("E4").Value = 9988273

# Row 5 (new row)
$ws.Range("A5").Value = "dev_1"
$ws.Range("B5").Value = "check_1233"
$ws.Range("C5").Value = "tester"
$ws.Range("D5").Value = "test@test.com"
$ws.Range("E5").Value = 9988273
